$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacity")

# --- Data value updates (wind 2025 FI etc.) ---
$ws.Range("F2").Value = 9600
$ws.Range("F4").Value = 3400
$ws.Range("F5").Value = 7600
$ws.Range("F7").Value = 2280

# --- Clear the stray empty-but-styled B cells in the Distributed Energy block ---
$ws.Range("B31:B34").ClearContents()

# --- Populate row 35 with a new Distributed Energy entry for FI00 / Onshore Wind ---
$ws.Range("A35").Value = "FI00"
$ws.Range("C35").Value = "Onshore Wind"
$ws.Range("D35").Value = "Distributed Energy"
$ws.Range("E35").Value = 2040
$ws.Range("F35").Value = 21000

# --- AutoFilter over the data range, with the matching hidden _FilterDatabase name ---
$ws.Range("A1:J35").AutoFilter()
$fdb = $ws.Names.Add("_FilterDatabase", "=Capacity!`$A`$1:`$J`$35")
$fdb.Visible = $false

# --- Update the saved view state (selection / active cell) ---
$ws.Activate()
$ws.Range("F28").Select()

Write-Output "done"
